# Update SwaadSutra_Daily_2026-01-19.xlsx
# - Order (row 2 of "Daily Orders") items/total changed.
# - "Summary" total revenue adjusted to match.
# - "Items Breakdown" recomputed per-item quantities/revenue, with a new
#   "Upma" row added.

$wb = $excel.ActiveWorkbook

# --- Daily Orders sheet: update row 2 (order #18) ---
$orders = $wb.Worksheets.Item("Daily Orders")
$orders.Range("F2").Value = "Pohe x2, Wheat Chapati x2, Upma x1"
$orders.Range("G2").Value = 120

# --- Summary sheet: update total revenue ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("G2").Value = 405

# --- Items Breakdown sheet: rewrite rows to reflect the new item mix ---
$items = $wb.Worksheets.Item("Items Breakdown")

# Row 2: Wheat Chapati quantity/revenue increase
$items.Range("A2").Value = "Wheat Chapati"
$items.Range("B2").Value = 15
$items.Range("C2").Value = 225

# Row 3: 1 Plate Bhaji (unchanged values, now listed before Pohe)
$items.Range("A3").Value = "1 Plate Bhaji"
$items.Range("B3").Value = 3
$items.Range("C3").Value = 90

# Row 4: Pohe quantity/revenue decrease
$items.Range("A4").Value = "Pohe"
$items.Range("B4").Value = 2
$items.Range("C4").Value = 60

# Row 5: new Upma row
$items.Range("A5").Value = "Upma"
$items.Range("B5").Value = 1
$items.Range("C5").Value = 30
